$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2607.1428
$ws.Range("I40").Value = 2607.1428
$ws.Range("K40").Value = 2607.1428
$ws.Range("M40").Value = -2432.1428
$ws.Range("H64").Value = 3966.6667
$ws.Range("I64").Value = 3750
$ws.Range("J64").Value = 4400
$ws.Range("K64").Value = 3750
$ws.Range("L64").Value = 4400
$ws.Range("M64").Value = -3502
$ws.Range("N64").Value = -4896
$ws.Range("H67").Value = 3966.6667
$ws.Range("I67").Value = 3750
$ws.Range("J67").Value = 4400
$ws.Range("K67").Value = 3750
$ws.Range("L67").Value = 4400
$ws.Range("M67").Value = -2892
$ws.Range("N67").Value = -6116
$ws.Range("H70").Value = 1499.8077
$ws.Range("J70").Value = 1521.0714
$ws.Range("L70").Value = 4563.2142
$ws.Range("N70").Value = -5103.2142
$ws.Range("H73").Value = 1499.8077
$ws.Range("J73").Value = 1521.0714
$ws.Range("L73").Value = 4563.2142
$ws.Range("N73").Value = -6435.2142
$ws.Range("H74").Value = 3999.5
$ws.Range("I74").Value = 3999.5
$ws.Range("K74").Value = 3999.5
$ws.Range("M74").Value = -3063.5
$ws.Range("H77").Value = 3999.5
$ws.Range("I77").Value = 3999.5
$ws.Range("K77").Value = 19997.5
$ws.Range("M77").Value = -15317.5
$ws.Range("H97").Value = 3999.4
$ws.Range("J97").Value = 3999.4
$ws.Range("L97").Value = 11998.2
$ws.Range("N97").Value = -12990.2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7439.7
$ws.Range("I32").Value = 5726.0527
$ws.Range("K32").Value = 5726.0527
$ws.Range("M32").Value = -5439.0527
$ws.Range("H74").Value = 433.33334
$ws.Range("I74").Value = 400
$ws.Range("K74").Value = 400
$ws.Range("M74").Value = 474
$ws.Range("H77").Value = 433.33334
$ws.Range("I77").Value = 400
$ws.Range("K77").Value = 2000
$ws.Range("M77").Value = 2368
$ws.Range("H88").Value = 1997.6666
$ws.Range("I88").Value = 1996
$ws.Range("K88").Value = 1996
$ws.Range("M88").Value = -1590
$ws.Range("H91").Value = 1997.6666
$ws.Range("I91").Value = 1996
$ws.Range("K91").Value = 1996
$ws.Range("M91").Value = -592
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4000
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("H89").Value = 4000
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("H105").Value = 0
$ws.Range("I105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("M105").ClearContents()
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 750
$ws.Range("I4").Value = 750
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 750
$ws.Range("L4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -638
$ws.Range("H7").Value = 36.636364
$ws.Range("I7").Value = 11.166667
$ws.Range("K7").Value = 11.166667
$ws.Range("M7").Value = 101.833333
$ws.Range("H39").Value = 4500
$ws.Range("I39").Value = 4500
$ws.Range("K39").Value = 4500
$ws.Range("M39").Value = -4109
$ws.Range("H49").Value = 4500
$ws.Range("I49").Value = 4500
$ws.Range("K49").Value = 4500
$ws.Range("M49").Value = -4318
$ws.Range("H62").Value = 4950
$ws.Range("I62").Value = 4900
$ws.Range("J62").Value = 5000
$ws.Range("K62").Value = 4900
$ws.Range("L62").Value = 5000
$ws.Range("M62").Value = -4276
$ws.Range("N62").Value = -6248
$ws.Range("H65").Value = 4950
$ws.Range("I65").Value = 4900
$ws.Range("J65").Value = 5000
$ws.Range("K65").Value = 24500
$ws.Range("L65").Value = 25000
$ws.Range("M65").Value = -21380
$ws.Range("N65").Value = -31240
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 150
$ws.Range("I7").Value = 100
$ws.Range("K7").Value = 300
$ws.Range("M7").Value = -188
$ws.Range("H25").Value = 500
$ws.Range("I25").Value = 500
$ws.Range("J25").Value = 500
$ws.Range("K25").Value = 1500
$ws.Range("L25").Value = 1500
$ws.Range("M25").Value = -1331
$ws.Range("N25").Value = -1838
$ws.Range("H30").Value = 500
$ws.Range("I30").Value = 500
$ws.Range("J30").Value = 500
$ws.Range("K30").Value = 1500
$ws.Range("L30").Value = 1500
$ws.Range("M30").Value = -1398
$ws.Range("N30").Value = -1704
$ws.Range("H38").Value = 8499.666999999999
$ws.Range("I38").Value = 8250
$ws.Range("J38").Value = 8999
$ws.Range("K38").Value = 24750
$ws.Range("L38").Value = 26997
$ws.Range("M38").Value = -24403
$ws.Range("N38").Value = -27691
$ws.Range("H48").Value = 4584
$ws.Range("I48").Value = 200
$ws.Range("J48").Value = 8968
$ws.Range("K48").Value = 600
$ws.Range("L48").Value = 26904
$ws.Range("M48").Value = -350
$ws.Range("N48").Value = -27404
$ws.Range("H97").Value = 625
$ws.Range("I97").Value = 550
$ws.Range("K97").Value = 1650
$ws.Range("M97").Value = -1154
$ws.Range("H138").Value = 1999.9333
$ws.Range("I138").Value = 1999
$ws.Range("K138").Value = 5997
$ws.Range("M138").Value = -857
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 20600
$ws.Range("I80").Value = 2900
$ws.Range("K80").Value = 2900
$ws.Range("M80").Value = -1902
$ws.Range("H83").Value = 20600
$ws.Range("I83").Value = 2900
$ws.Range("K83").Value = 14500
$ws.Range("M83").Value = -9508
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 36667
$ws.Range("J2").Value = 36667
$ws.Range("L2").Value = 36667
$ws.Range("N2").Value = -36891
$ws.Range("H22").Value = 831.6667
$ws.Range("I22").Value = 495
$ws.Range("K22").Value = 495
$ws.Range("M22").Value = -200
$ws.Range("H27").Value = 831.6667
$ws.Range("I27").Value = 495
$ws.Range("K27").Value = 495
$ws.Range("M27").Value = -388
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 20000000
$ws.Range("J5").Value = 20000000
$ws.Range("L5").Value = 20000000
$ws.Range("N5").Value = -20000224
$ws.Range("H22").Value = 3890
$ws.Range("J22").Value = 3890
$ws.Range("L22").Value = 3890
$ws.Range("N22").Value = -4476
$ws.Range("H96").Value = 3095.8
$ws.Range("I96").Value = 3000
$ws.Range("K96").Value = 3000
$ws.Range("M96").Value = -1627
$ws.Range("H126").Value = 539.8
$ws.Range("I126").Value = 539.8
$ws.Range("K126").Value = 1619.4
$ws.Range("M126").Value = 850.6000000000001
